$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19 mirrors the formatting of row 18 (same column layout: TCID, Jira id,
# Description, Runmode, Results). Copy formats from row 18 into row 19 first so
# borders / wrap-text / etc. match the rest of the table.
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match row 18's row height (45pt), used for the wrapped multi-line entries.
$ws.Rows("19").RowHeight = 45

# Shared strings get appended in the same order the new values are written below:
# Jira id (B) -> Description (C) -> TCID (A), matching the sharedStrings.xml diff.
$ws.Range("B19").Value = "OPQA-4554||OPQA-4555||OPQA-4553||OPQA-4541"
$ws.Range("C19").Value = "Verify that error message ""Invalid email/password. Please try again."" should be displayed when user enters incorrect password.||Verify that error message should be displayed when user provides correct STeAM credentials but the STeAM account is not entitled to DRA\IPA application. As per wireframe||Verify that the DRA\IPA application overview page shall be able to be opened in a separate browser window on demand. (e.g. from a link on the Step Up Authentication Modal)|| user be should be presented the appropriate ""not entitled"" modal as an error overlay on the target application sign in page."
$ws.Range("A19").Value = "IPAIAM0058"

# Runmode (D19) stays "Y", same as the rest of the table; Results (E19) stays blank.
$ws.Range("D19").Value = "Y"

# Update the active selection as recorded in the saved workbook.
$ws.Range("C21").Select() | Out-Null
